$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot to the latest scrape.
# Cells whose new text would otherwise auto-parse as a number are forced
# to remain plain text (matching the source data, which stores everything
# as inline strings) and then restored to the default "Normal" style so no
# stray number-format is left behind.
$ws.Range('D2').Value = '26.992.68'
$ws.Range('E2').Value = '  +2.78%  '
$ws.Range('D3').Value = '1.652.34'
$ws.Range('E3').Value = '  +3.62%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +1.62%  '
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.59'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0864'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').Value = '1.886.01'
$ws.Range('E12').Value = '  +3.72%  '
$ws.Range('D13').Value = '1.645.06'
$ws.Range('E13').Value = '  +3.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.519'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('D17').Value = '27.030.42'
$ws.Range('E17').Value = '  +3.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '237.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.61%  '
$ws.Range('E19').Value = '  +2.58%  '
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.22'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  +2.28%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.18%  '
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('E31').Value = '  +1.56%  '
$ws.Range('D32').Value = '1.528.91'
$ws.Range('E32').Value = '  +4.40%  '
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('E34').Value = '  +2.98%  '
$ws.Range('E35').Value = '  +8.26%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.884'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.90%  '
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.93'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.66%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.26'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.98%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '66.19'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +9.50%  '
$ws.Range('D44').Value = '1.793.14'
$ws.Range('E44').Value = '  +3.56%  '
$ws.Range('E45').Value = '  +2.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.921'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.14'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.96%  '
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0503'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('E51').Value = '  +3.00%  '
